# Add "DICOM:SeriesNumber" as a new default DICOM metadata column.
#
# The workbook has a "Files" sheet whose row 1 holds column headers,
# several of which are DICOM metadata fields (columns J..O):
#   J: DICOM:Manufacturer
#   K: DICOM:ManufacturerModelName
#   L: DICOM:Modality
#   M: DICOM:StationName
#   N: DICOM:StudyDate
#   O: DICOM:StudyDescription
#   P: SubjectName
#   Q: ID
#   R: Directory
#   S: filename
#
# We insert a new "DICOM:SeriesNumber" column right after Modality (so
# before the old StationName column), which pushes the existing M..S
# columns one place to the right (becoming N..T), and fill in the
# SeriesNumber value (6168) for every data row that already carries
# DICOM metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at M; existing M:S shift right to N:T.
$ws.Columns("M:M").Insert()

# New header cell for the inserted column.
$ws.Range("M1").Value = "DICOM:SeriesNumber"

# Rows that already have DICOM metadata (scanned file rows) get the
# SeriesNumber value too.
$dataRows = @(3, 4, 5, 6, 7, 10, 12, 14, 16)
foreach ($r in $dataRows) {
    $ws.Range("M$r").Value = 6168
}

# Reflect the selection used when the column was added (whole column M).
$ws.Range("M1:M1048576").Select()
